# Daily attendance processing - normalize "Recorded By" (column G) entries
# so that the literal audit-trail tag "System" (exact case) always appears
# first in the comma-separated list of recorders, swapping it with
# whichever entry was first. Mirrors the upstream "system" account being
# re-flagged ahead of the collaborator/backup addresses that recorded
# each session.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $v = $cell.Value2

    if ($v -ne $null -and $v -ne "") {
        $parts = $v.Split(",")

        if ($parts.Length -ge 2) {
            $trimmed = @()
            foreach ($p in $parts) {
                $trimmed += $p.Trim()
            }

            $lastIdx = $trimmed.Length - 1

            # Only act when the trailing entry is exactly "System" -
            # leading "System" entries (e.g. "System, admin@admin.com")
            # and lists with no "System" tag are left untouched.
            if ($trimmed[$lastIdx] -eq "System") {
                $newParts = $trimmed
                $tmp = $newParts[0]
                $newParts[0] = $newParts[$lastIdx]
                $newParts[$lastIdx] = $tmp

                $newValue = [string]::Join(", ", $newParts)
                $cell.Value2 = $newValue
            }
        }
    }
}
